# Converts a target EMU value to a point value that, once PowerPoint's COM
# layer stores it as a single-precision float and re-derives EMU from it
# (points * 12700, truncated), reproduces exactly that EMU value. This
# compensates for float32 rounding/truncation in the Shape.Top/.Left setters.
function ConvertTo-ComPoints($targetEmu) {
    $base = $targetEmu / 12700.0
    for ($i = 0; $i -le 4000; $i++) {
        $candidate = $base + ($i * 0.0000005)
        $f = [single]$candidate
        $emu = [int]($f * 12700)
        if ($emu -eq $targetEmu) {
            return $candidate
        }
    }
    return $base
}

$p = $ppt.ActivePresentation

# --- Slide 1: subtitle text tweak ---------------------------------------
# "MSE800 - Professional Software Engineering" -> "MSE800 (Professional Software Engineering)"
$s1 = $p.Slides.Item(1)
$subtitle = $s1.Shapes.Item(2)
$para2 = $subtitle.TextFrame.TextRange.Paragraphs(2)
$run = $para2.Runs(1)
$run.Text = "MSE800 (Professional Software Engineering)"

# --- Slide 21 ("Thank You..."): nudge title box down --------------------
$s21 = $p.Slides.Item($p.Slides.Count)
$title = $s21.Shapes.Item(1)
$title.Top = ConvertTo-ComPoints 1908668
